$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 426 ("「仮にこれが批判の数としよう」" post), shifting
# all subsequent rows up by one.
$ws.Rows.Item(426).Delete()
